$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "浙江杭州滨江中南乐游城店_会员_借记") {
        $ws.Name = "浙江杭州滨江中南乐游城店_借记"
    }
    elseif ($ws.Name -eq "浙江杭州三墩地铁站店_会员_借记") {
        $ws.Name = "浙江杭州三墩地铁站店_借记"
    }
}
